# Logged Week 15 and simulated Week 16
# Update the "R" row (row 3) target-depth totals on both the OFF and DEF
# sheets to reflect the newly logged Week 15 data plus the simulated
# Week 16 data.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 409
$wsOff.Range("C3").Value = 299
$wsOff.Range("D3").Value = 118
$wsOff.Range("E3").Value = 50
$wsOff.Range("G3").Value = 5

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 467
$wsDef.Range("C3").Value = 326
$wsDef.Range("D3").Value = 99
$wsDef.Range("E3").Value = 48
$wsDef.Range("G3").Value = 5
